# Insert a new data row before the current row 388 ("2022-10-21" / Brasil / Primera / 228)
# by duplicating the row that will end up at 389 (keeps all styles/constant columns intact),
# then overwrite the cells that differ for this new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(388).Copy()
$ws.Rows.Item(388).Insert()

# Fill in the values specific to the newly inserted row 388
$ws.Range("D388").Value2 = 44855
$ws.Range("L388").Value = "Primera"
$ws.Range("M388").Value = 228
$ws.Range("N388").Value = 8000
$ws.Range("O388").Value = 8000
$ws.Range("P388").Value = 8000
$ws.Range("R388").Value = "Brasil"
$ws.Range("S388").Value = 2000
